$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update the query text in B2 (drop the trailing Cohort coalesce line) ---
$newB2 = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
  WHERE labels(parent)[0] IN ["sample"]
OPTIONAL MATCH (f)-[*]->(c:case)<--(demo:demographic)
OPTIONAL MATCH (s:study)<-[*]-(c)
OPTIONAL MATCH (c)<--(diag:diagnosis)
OPTIONAL  MATCH (samp:sample)-->(c)
OPTIONAL  MATCH (co:cohort)<-[*]-(c)
  WITH DISTINCT c, s, demo, diag, co
RETURN coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``
"@
$ws.Range("B2").Value = $newB2

# --- View / selection changes ---
$ws.Range("B2").Select()
$excel.ActiveWindow.Zoom = 85

# --- Row height adjustments (rows 2-4 shrink slightly to match new render) ---
$ws.Rows.Item(2).RowHeight = 288
$ws.Rows.Item(3).RowHeight = 288
$ws.Rows.Item(4).RowHeight = 259.2
